$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 6.795
$ws.Range("B13").Value = 6.606
$ws.Range("B16").Value = 5.787000000000001
$ws.Range("B18").Value = 6.313000000000001
$ws.Range("B20").Value = 6.661
